# Updating scenario file names (NOT FINAL VERSIONS)
#
# The "Platform Coverage" sheet's MDA age-bracket rows were out of sync
# with their coverage percentages. This inserts a new row for the
# "2-15" age bracket (which keeps the high 0.8 coverage values that used
# to incorrectly live on the "min age" row), shifts the min/max ages of
# the other MDA rows down one bracket, and adds a new 50-65 bracket row
# with 0.5 coverage - mirroring the already-present 15-50 bracket.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

$coverageCols = @("P", "R", "T", "V", "X", "Z", "AB", "AD")

# Insert a new row at position 3; this pushes the existing rows 3-7
# down to rows 4-8 (EPI / School / Out-of-school campaign rows all
# move down untouched).
$ws.Rows.Item(3).Insert()

# New row 3 becomes the "2-15" MDA bracket, carrying the 0.8 coverage
# values that used to sit (incorrectly) on row 2.
$ws.Range("A3").Value = "All"
$ws.Range("B3").Value = "Treatment"
$ws.Range("C3").Value = "Campaign"
$ws.Range("D3").Value = "MDA"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 15
foreach ($col in $coverageCols) {
    $ws.Range($col + "3").Value = 0.8
}

# Row 2 ("min age" bracket) now covers ages 5-15 and no longer carries
# the later-year coverage figures (those moved to row 3 above).
$ws.Range("F2").Value = 5
foreach ($col in $coverageCols) {
    $ws.Range($col + "2").ClearContents()
}

$wb.Save()
